$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
    }
}

# 1. Merge the paragraph break between "...up to. " and "Although this led to higher..."
Replace-Text "assembly line was up to. `rAlthough this led to higher " "assembly line was up to. Although this led to higher "

# 2. Remove the _GoBack bookmark before "Smith is pointing out a"
$bm = $d.Bookmarks("_GoBack")
$bm.Delete() | Out-Null

# 3. "As in a socialist economy" -> "As in a planned, socialist economy"
Replace-Text "As in a socialist economy" "As in a planned, socialist economy"

# 4. Merge "knowledge commissars. " paragraph break into text, but insert "major " and re-split the
#    paragraph at the same place ("A " | "major" | " problem...") - content stays split into two paragraphs
#    as before: only the text of the second paragraph changes from "A problem" to "A major problem".
Replace-Text "A problem with this approach is that as products become more complicated" "A major problem with this approach is that as products become more complicated"

# 5. "...the evolving understanding of all of the other knowledge workers involved in the product: ..."
#    -> "...the evolving understanding of all of the other knowledge workers producing the product: ..."
Replace-Text "other knowledge workers involved in the product" "other knowledge workers producing the product"

# 6. Insert italic "almost " before "everyone" ("...typically encompasses everyone working...")
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("everyone working on the product.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertRng = $d.Range($rng.Start, $rng.Start)
    $insertRng.InsertBefore("almost ")
    $insertRng2 = $d.Range($rng.Start, $rng.Start + 7)
    $insertRng2.Font.Italic = 1
} else {
    Write-Host "NOT FOUND: everyone working on the product."
}
$found = $null

# 7. "production follow from the nature of knowledge workers cooperating" ->
#    "production follow from recognizing these realities concerning knowledge workers cooperating"
Replace-Text "production follow from the nature of knowledge workers cooperating" "production follow from recognizing these realities concerning knowledge workers cooperating"

# 8. "...those small batches, " -> "...those small batches; " (comma -> semicolon)
Replace-Text "allowing end users to comment on the work done in those small batches, " "allowing end users to comment on the work done in those small batches; "

# 9. "A rigid division of labor hinders" -> "Given the above realities, a rigid division of labor hinders"
Replace-Text "A rigid division of labor hinders" "Given the above realities, a rigid division of labor hinders"

# 10. "confined to a narrow silos based" -> "confined to narrow silos based" (also drops proofErr tags)
Replace-Text "confined to a narrow silos based" "confined to narrow silos based"

# 11. Insert the _GoBack bookmark right after "If workers are confined to "
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("If workers are confined to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $collapsed = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $collapsed) | Out-Null
} else {
    Write-Host "NOT FOUND: If workers are confined to "
}
